# Variant 8 (ЕГЭ, задание 18) — "путь минимальной стоимости" проверка ответов.
# Adds the worked-solution rows (15-25) below the existing 11x12 data grid:
#   - Row 15: best cost reachable from row 1, moving right or down into row 16
#   - Rows 16-24: DP table seeded from rows 2-10, each cell = min(move right, move down)
#   - Row 25: plain running sum seeded from row 11 (no branching)
# and restores the view/print state that Excel re-serialises on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 15 (built from row 1) ----
$ws.Range("A15").Formula      = '=MIN(SUM(A1, B15), SUM(A1, $A16))'
$ws.Range("B15:K15").Formula  = '=MIN(SUM(B1, C15), SUM(B1, $A16))'
$ws.Range("L15").Formula      = '=SUM(L1, $A16)'

# ---- Rows 16-24 (built from rows 2-10), one shared formula per column ----
$ws.Range("A16:A24").Formula = '=MIN(SUM(A2, B16), SUM(A2, $A17))'
$ws.Range("B16:B24").Formula = '=MIN(SUM(B2, C16), SUM(B2, $A17))'
$ws.Range("C16:C24").Formula = '=MIN(SUM(C2, D16), SUM(C2, $A17))'
$ws.Range("D16:D24").Formula = '=MIN(SUM(D2, E16), SUM(D2, $A17))'
$ws.Range("E16:E24").Formula = '=MIN(SUM(E2, F16), SUM(E2, $A17))'
$ws.Range("F16:F24").Formula = '=MIN(SUM(F2, G16), SUM(F2, $A17))'
$ws.Range("G16:G24").Formula = '=MIN(SUM(G2, H16), SUM(G2, $A17))'
$ws.Range("H16:H24").Formula = '=MIN(SUM(H2, I16), SUM(H2, $A17))'
$ws.Range("I16:I24").Formula = '=MIN(SUM(I2, J16), SUM(I2, $A17))'
$ws.Range("J16:J24").Formula = '=MIN(SUM(J2, K16), SUM(J2, $A17))'
$ws.Range("K16:K24").Formula = '=MIN(SUM(K2, L16), SUM(K2, $A17))'
$ws.Range("L16:L24").Formula = '=SUM(L2, $A17)'

# ---- Row 25 (built from row 11, plain chained SUM) ----
$ws.Range("A25").Formula     = '=SUM(A11, B25)'
$ws.Range("B25:L25").Formula = '=SUM(B11, C25)'

# ---- View state: active cell moves to S13 ----
[void]$ws.Range("S13").Select()

# ---- Print setup: A4 portrait (paperSize 9 = A4) ----
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
